$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = -0.2705469131469726
$ws.Range("B2").Value = 0.4402385950088501
$ws.Range("C2").Value = -1.731513738632202
$ws.Range("A3").Value = -0.6277971267700195
$ws.Range("B3").Value = 0.0076048374176025
$ws.Range("C3").Value = -1.507715225219727
$ws.Range("A4").Value = -2.914698600769043
$ws.Range("B4").Value = -1.449564576148987
$ws.Range("C4").Value = -3.32840347290039
$ws.Range("A5").Value = 1.028462886810303
$ws.Range("B5").Value = -0.5502710342407227
$ws.Range("C5").Value = -4.842555046081543
$ws.Range("A6").Value = -1.341280460357666
$ws.Range("B6").Value = -2.225003957748413
$ws.Range("C6").Value = -6.344600677490234
$ws.Range("A7").Value = 8.574896812438965
$ws.Range("B7").Value = 0.6133027076721191
$ws.Range("C7").Value = -6.888121604919434
$ws.Range("A8").Value = -6.096681118011475
$ws.Range("B8").Value = 0.8472604751586914
$ws.Range("C8").Value = 14.72706890106201
$ws.Range("A9").Value = 4.274323463439941
$ws.Range("B9").Value = -4.468049049377441
$ws.Range("C9").Value = -6.856836795806885
$ws.Range("A10").Value = -4.518700122833252
$ws.Range("B10").Value = -1.648021101951599
$ws.Range("C10").Value = -0.9248533248901368
$ws.Range("A11").Value = 9.755411148071287
$ws.Range("B11").Value = 3.367114305496216
$ws.Range("C11").Value = 2.822277307510376
$ws.Range("A12").Value = 1.561064720153809
$ws.Range("B12").Value = 0.1129603385925293
$ws.Range("C12").Value = -0.9029455184936525
$ws.Range("A13").Value = 5.92741584777832
$ws.Range("B13").Value = -0.8555939197540283
$ws.Range("C13").Value = 4.797466278076172
$ws.Range("A14").Value = 1.122594833374023
$ws.Range("B14").Value = 1.295500755310059
$ws.Range("C14").Value = -1.442571401596069
$ws.Range("A15").Value = 0.5986118316650391
$ws.Range("B15").Value = 0.4096674025058746
$ws.Range("C15").Value = -0.6679027080535889
$ws.Range("A16").Value = 0.0388402938842773
$ws.Range("B16").Value = 0.3524296283721924
$ws.Range("C16").Value = -1.101761341094971
$ws.Range("A17").Value = -0.1728830337524414
$ws.Range("B17").Value = 0.6193998456001282
$ws.Range("C17").Value = -0.6873818635940552
$ws.Range("A18").Value = 0.4876585006713867
$ws.Range("B18").Value = 0.6636635065078735
$ws.Range("C18").Value = -0.9166454076766968
$ws.Range("A19").Value = -0.1092472076416015
$ws.Range("B19").Value = 0.732629120349884
$ws.Range("C19").Value = -1.016466021537781
$ws.Range("A20").Value = 0.4153709411621094
$ws.Range("B20").Value = 0.5096800327301025
$ws.Range("C20").Value = -0.7671611309051514
$ws.Range("A21").Value = 0.17730712890625
$ws.Range("B21").Value = 0.6253083348274231
$ws.Range("C21").Value = -0.8837988376617432
